# Adds 11 new match rows (rows 51-61, index 50-60) to the Kuwait Premier
# League 2023-2024 sheet, mirroring the commit "Atualizado por script em
# 02-01-2024 20:45". Columns: Indice, pais, torneio, temporada, data_partida,
# home, home_ft_gols, away, away_ft_gols, home_opening_odds,
# home_opening_data_hora, home_closing_odds, home_closing_data_hora,
# draw_opening_odds, draw_opening_data_hora, draw_closing_odds,
# draw_closing_data_hora, away_opening_odds, away_opening_data_hora,
# away_closing_odds, away_closing_data_hora, url_partida.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One line per new row, tab-separated, columns A..V in sheet order.
$data = @"
50	kuwait	premier-league	2023-2024	45283.64583333334	Al Qadisiya	2	Al Jahra	1	1.32	23/12/2023 03:42	1.28	23/12/2023 15:26	4.92	23/12/2023 03:42	5.28	23/12/2023 15:26	7.54	23/12/2023 03:42	8.99	23/12/2023 15:26	https://www.betexplorer.com/football/kuwait/premier-league/al-qadisiya-al-jahra/UZMFG1vt/
51	kuwait	premier-league	2023-2024	45284.64583333334	Al-Fahaheel	1	Khaitan	0	2.05	24/12/2023 03:41	1.75	24/12/2023 15:23	3.42	24/12/2023 03:41	3.67	24/12/2023 15:23	3.21	24/12/2023 03:41	4.17	24/12/2023 15:23	https://www.betexplorer.com/football/kuwait/premier-league/al-fahaheel-khaitan/zXPNEu9h/
52	kuwait	premier-league	2023-2024	45284.75	Al Kuwait	6	Kazma SC	0	1.5	24/12/2023 06:12	1.4	24/12/2023 17:57	4.26	24/12/2023 06:12	4.87	24/12/2023 17:59	5.32	24/12/2023 06:12	6.14	24/12/2023 17:59	https://www.betexplorer.com/football/kuwait/premier-league/al-kuwait-kazma-sc/UeK784nJ/
53	kuwait	premier-league	2023-2024	45285.64583333334	Al Arabi	1	Al Shabab	0	1.3	25/12/2023 03:42	1.24	25/12/2023 15:13	5.05	25/12/2023 03:42	5.73	25/12/2023 15:28	8.09	25/12/2023 03:42	9.31	25/12/2023 15:28	https://www.betexplorer.com/football/kuwait/premier-league/al-arabi-kuwait-al-shabab/tAQJFLgn/
54	kuwait	premier-league	2023-2024	45285.75	Al Naser	1	Al Salmiya	1	2.24	25/12/2023 06:12	2.15	25/12/2023 17:53	3.43	25/12/2023 06:12	3.63	25/12/2023 17:53	2.83	25/12/2023 06:12	2.9	25/12/2023 17:53	https://www.betexplorer.com/football/kuwait/premier-league/al-naser-al-salmiya/j5L39pXC/
55	kuwait	premier-league	2023-2024	45289.59027777778	Al Qadisiya	2	Al-Fahaheel	1	1.41	28/12/2023 13:08	1.45	29/12/2023 14:03	4.04	28/12/2023 13:08	4.28	29/12/2023 14:03	6.15	28/12/2023 13:08	6.15	29/12/2023 14:03	https://www.betexplorer.com/football/kuwait/premier-league/al-qadisiya-al-fahaheel/ENORDaOb/
56	kuwait	premier-league	2023-2024	45289.70138888889	Khaitan	0	Al Arabi	4	7.53	28/12/2023 13:08	8.98	29/12/2023 16:46	5.16	28/12/2023 13:08	5.62	29/12/2023 16:46	1.25	28/12/2023 13:08	1.25	29/12/2023 16:40	https://www.betexplorer.com/football/kuwait/premier-league/khaitan-al-arabi-kuwait/WbJWCJw5/
57	kuwait	premier-league	2023-2024	45290.59027777778	Al Shabab	0	Al Naser	1	3.35	30/12/2023 02:12	4.1	30/12/2023 12:15	3.53	30/12/2023 02:12	3.72	30/12/2023 12:15	1.95	30/12/2023 02:12	1.75	30/12/2023 12:15	https://www.betexplorer.com/football/kuwait/premier-league/al-shabab-al-naser/AkIzCwgB/
58	kuwait	premier-league	2023-2024	45290.70138888889	Al Salmiya	2	Al Kuwait	2	5.38	30/12/2023 05:12	4.46	30/12/2023 16:47	4.29	30/12/2023 05:12	4.16	30/12/2023 16:47	1.49	30/12/2023 05:12	1.61	30/12/2023 16:47	https://www.betexplorer.com/football/kuwait/premier-league/al-salmiya-al-kuwait/IFLrAHNN/
59	kuwait	premier-league	2023-2024	45291.64583333334	Al Jahra	0	Kazma SC	2	3.57	30/12/2023 15:13	3.5	31/12/2023 15:28	3.51	30/12/2023 15:13	3.75	31/12/2023 15:28	1.81	30/12/2023 15:13	1.88	31/12/2023 15:28	https://www.betexplorer.com/football/kuwait/premier-league/al-jahra-kazma-sc/4CHvBc8H/
60	kuwait	premier-league	2023-2024	45293.64930555555	Al Arabi	2	Al Qadisiya	2	2.17	02/01/2024 03:42	2.1	02/01/2024 15:34	3.35	02/01/2024 03:42	3.21	02/01/2024 15:33	2.98	02/01/2024 03:42	3.36	02/01/2024 15:34	https://www.betexplorer.com/football/kuwait/premier-league/al-arabi-kuwait-al-qadisiya/0vBm9ywU/
"@

# Column kind per letter A..V: "n" = numeric, "s" = text.
$colKinds = @('n','s','s','s','n','s','n','s','n','n','s','n','s','n','s','n','s','n','s','n','s','s')

$lines = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

$startRow = 51
$rowIdx = $startRow

# Template cells whose formatting (borders/font/number-format) must be
# reproduced on the new rows without inventing new style records: row 2
# already carries the "index" style (column A) and the "match datetime"
# style (column E) used throughout the sheet.
$styleSourceA = $ws.Cells.Item(2, 1)
$styleSourceE = $ws.Cells.Item(2, 5)

foreach ($line in $lines) {
    $fields = $line -split "`t"

    for ($col = 1; $col -le 22; $col++) {
        $cell = $ws.Cells.Item($rowIdx, $col)
        $raw = $fields[$col - 1]
        $kind = $colKinds[$col - 1]

        if ($col -eq 1) {
            $styleSourceA.Copy()
            $cell.PasteSpecial(-4122)
        } elseif ($col -eq 5) {
            $styleSourceE.Copy()
            $cell.PasteSpecial(-4122)
        }

        if ($kind -eq 'n') {
            $cell.Value = [double]$raw
        } else {
            $cell.Value = $raw
        }
    }

    $rowIdx++
}

$excel.CutCopyMode = 0
